# Edits writeup.docx per the commit: "made some slight changes"

$d = $word.ActiveDocument

# --- 1) Paragraph 6: wrap the word "presented" in a bookmark (mirrors a
#        Word auto-generated __DdeLink bookmark left by a paste/DDE edit). ---
$rng = $d.Content
$rng.Find.Execute("presented", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("__DdeLink__40_1327815132", $rng)

# --- 2) Paragraph 7: "as search problem" -> "as a search problem" ---
$rng = $d.Content
$rng.Find.Execute("as search problem", $true, $false, $false, $false, $false, $true, 1, $false, "as a search problem", 2)

# --- 3) Paragraph 9: "These inputs come for a predefined" -> "...come from a predefined" ---
$rng = $d.Content
$rng.Find.Execute("inputs come for a predefined", $true, $false, $false, $false, $false, $true, 1, $false, "inputs come from a predefined", 2)

# --- 4) Paragraph 12: possible states count/formula update ---
$rng = $d.Content
$rng.Find.Execute("The total number of possible states is 9! = 362,880", $true, $false, $false, $false, $false, $true, 1, $false, "The total number of SOLVABLE possible states is 9!/2 = 181,440", 2)

# --- 5) Paragraph 16: drop the parenthetical and add the 3x3-matrix sentence ---
$rng = $d.Content
$rng.Find.Execute("The possible actions are up, down, left, and right while being limited to remain within the tiled board (ie in range [0:2][0:2]). ", $true, $false, $false, $false, $false, $true, 1, $false, "The possible actions are up, down, left, and right while being limited to remain within the tiled board. This means remaing within a 3X3 matrix. ", 2)

# --- 6) Insert two new bullets ("Successor State" / its sub-bullet) right
#        after "Right: swap tiles [E] and [E+1]" and before "Goal:" ---
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "Right: swap tiles [E] and [E+1]") {
        $targetIdx = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIdx)
$p.Range.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Item($targetIdx + 1)
$newPara1.Range.ListFormat.ListLevelNumber = 1
$newPara1.Range.Text = "Successor State"

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($targetIdx + 2)
$newPara2.Range.ListFormat.ListLevelNumber = 2
$newPara2.Range.Text = "Position based on transition model gives the successor state"

# --- 7) Remove the stray _GoBack bookmark at the end of the document ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
